$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.197.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").Value = "'3.402.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.01%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'584.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = "'182.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.69%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("D9").Value = "'0.203"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +10.76%  '
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("D11").Value = "'48.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("D12").Value = "'0.0000289"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.52%  '
$ws.Range("D13").Value = "'689.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = "'8.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.12%  '
$ws.Range("D15").Value = "'3.951.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").Value = "'70.120.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").Value = "'3.404.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = "'17.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = "'11.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("E21").Value = '  +2.47%  '
$ws.Range("D22").Value = "'17.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.75%  '
$ws.Range("D23").Value = "'5.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = "'102.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("D25").Value = "'3.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("D27").Value = "'9.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.32%  '
$ws.Range("D28").Value = "'33.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.35%  '
$ws.Range("E29").Value = '  +3.36%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").Value = "'3.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.94%  '
$ws.Range("D32").Value = "'11.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("D33").Value = "'557.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("D35").Value = "'58.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = "'3.653.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("E38").Value = '  +3.63%  '
$ws.Range("D39").Value = "'35.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("E40").Value = '  +8.98%  '
$ws.Range("E41").Value = '  +5.34%  '
$ws.Range("D42").Value = "'2.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.10%  '
$ws.Range("D43").Value = "'0.0431"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.76%  '
$ws.Range("D44").Value = "'0.341"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.80%  '
$ws.Range("E45").Value = '  +2.41%  '
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E47").Value = '  +4.44%  '
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").Value = "'130.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("E50").Value = '  -0.01%  '
